{"js": "// Ran test cases on contact form.\n//\n// Updates the sprint-plan bullet list:\n//   Sprint 2 - append \"+ set up test suite\"\n//   Sprint 3 - prepend \"DB Implementation and hosting + REST backend + \" to \"Message Service\"\n//   Sprint 4 - append \" (Chat) + try unit testing\" after \"Service\"\n//   Sprint 5 - \"Window Handler and Board Menu\" -> \"Implementing socket.io\"\n//   Sprint 6 - rewrite tail + move the \"_GoBack\" bookmark here\n//   Sprint 7 - \"Navigation and design\" -> \"Navigation, Backbone History, design (mobile first) + testing\"\n//   Sprint 8 - \"Navigation and design\" (+ old bookmark) -> \"Remaining Bug fixes + design (mobile first) + testing\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Replace the single occurrence of `target` inside paragraph `index` with `replacement`,\n// scoping the search to that paragraph so identical phrases elsewhere are untouched.\nasync function replaceInParagraph(index, target, replacement) {\n  const paragraph = paragraphs.items[index];\n  const found = paragraph.search(target, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${target}\" in paragraph ${index}, found ${found.items.length}`\n    );\n  }\n  found.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// Sprint 2: \"... DB Implementation and hosting.\" -> \"... DB Implementation and hosting + set up test suite.\"\nawait replaceInParagraph(\n  0,\n  \"DB Implementation and hosting.\",\n  \"DB Implementation and hosting + set up test suite.\"\n);\n\n// Sprint 3: \"... Message Service.\" -> \"... DB Implementation and hosting + REST backend + Message Service.\"\nawait replaceInParagraph(\n  1,\n  \"Message Service.\",\n  \"DB Implementation and hosting + REST backend + Message Service.\"\n);\n\n// Sprint 4: \"... Message Service.\" -> \"... Message Service (Chat) + try unit testing.\"\nawait replaceInParagraph(\n  2,\n  \"Message Service.\",\n  \"Message Service (Chat) + try unit testing.\"\n);\n\n// Sprint 5: \"Window Handler and Board Menu\" -> \"Implementing socket.io\"\nawait replaceInParagraph(3, \"Window Handler and Board Menu\", \"Implementing socket.io\");\n\n// Sprint 6: rewrite the tail of the bullet. A temporary marker (U+2603) stands in for the\n// \"_GoBack\" bookmark position so it can be located precisely and replaced below.\nawait replaceInParagraph(\n  4,\n  \" Application Nav and general information about me + contact form.\",\n  \" Application Nav, design (mobile first) \\u2603+ contact form + testing.\"\n);\n\n// Sprint 7: \"Navigation and design\" -> \"Navigation, Backbone History, design (mobile first) + testing\"\nawait replaceInParagraph(\n  5,\n  \"Navigation and design\",\n  \"Navigation, Backbone History, design (mobile first) + testing\"\n);\n\n// Sprint 8: \"Navigation and design\" -> \"Remaining Bug fixes + design (mobile first) + testing\".\n// The \"_GoBack\" bookmark used to live at the end of this run; delete it here (it is recreated\n// at its new home in Sprint 6 immediately below).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait replaceInParagraph(\n  6,\n  \"Navigation and design\",\n  \"Remaining Bug fixes + design (mobile first) + testing\"\n);\n\n// Drop the \"_GoBack\" bookmark at the marker left in Sprint 6, then remove the marker itself.\n{\n  const paragraph = paragraphs.items[4];\n  const marker = paragraph.search(\"\\u2603\", { matchCase: true });\n  marker.load(\"text\");\n  await context.sync();\n  if (marker.items.length !== 1) {\n    throw new Error(`Expected exactly one bookmark marker, found ${marker.items.length}`);\n  }\n  const markerRange = marker.items[0];\n  markerRange.insertText(\"\", \"Replace\");\n  await context.sync();\n  const collapsed = markerRange.getRange(\"Start\");\n  collapsed.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Ran test cases on contact form.\n#\n# Updates the sprint-plan bullet list:\n#   Sprint 2 - append \"+ set up test suite\"\n#   Sprint 3 - prepend \"DB Implementation and hosting + REST backend + \" to \"Message Service\"\n#   Sprint 4 - append \" (Chat) + try unit testing\" after \"Service\"\n#   Sprint 5 - \"Window Handler and Board Menu\" -> \"Implementing socket.io\"\n#   Sprint 6 - rewrite tail + move the \"_GoBack\" bookmark here\n#   Sprint 7 - \"Navigation and design\" -> \"Navigation, Backbone History, design (mobile first) + testing\"\n#   Sprint 8 - \"Navigation and design\" (+ old bookmark) -> \"Remaining Bug fixes + design (mobile first) + testing\"\n\n$d = $word.ActiveDocument\n\n# Placeholder character used to mark, temporarily, where the \"_GoBack\" bookmark should end up.\n$marker = [char]0x2603\n\n# Replaces the single occurrence of $Target inside paragraph number $Index (1-based) with\n# $Replacement. Scoping Find to that paragraph's Range keeps identical phrases elsewhere untouched.\nfunction Replace-InParagraph($Index, $Target, $Replacement) {\n    $rng = $d.Paragraphs($Index).Range\n    $found = $rng.Find.Execute($Target)\n    if (-not $found) {\n        throw (\"Could not find '\" + $Target + \"' in paragraph \" + $Index)\n    }\n    $rng.Text = $Replacement\n}\n\n# Sprint 2: \"... DB Implementation and hosting.\" -> \"... DB Implementation and hosting + set up test suite.\"\nReplace-InParagraph 1 \"DB Implementation and hosting.\" \"DB Implementation and hosting + set up test suite.\"\n\n# Sprint 3: \"... Message Service.\" -> \"... DB Implementation and hosting + REST backend + Message Service.\"\nReplace-InParagraph 2 \"Message Service.\" \"DB Implementation and hosting + REST backend + Message Service.\"\n\n# Sprint 4: \"... Message Service.\" -> \"... Message Service (Chat) + try unit testing.\"\nReplace-InParagraph 3 \"Message Service.\" \"Message Service (Chat) + try unit testing.\"\n\n# Sprint 5: \"Window Handler and Board Menu\" -> \"Implementing socket.io\"\nReplace-InParagraph 4 \"Window Handler and Board Menu\" \"Implementing socket.io\"\n\n# Sprint 6: rewrite the tail of the bullet, dropping in the marker where the bookmark will go.\nReplace-InParagraph 5 \"Application Nav and general information about me + contact form.\" (\"Application Nav, design (mobile first) \" + $marker + \"+ contact form + testing.\")\n\n# Sprint 8: \"Navigation and design\" -> \"Remaining Bug fixes + design (mobile first) + testing\".\n# The \"_GoBack\" bookmark used to sit at the end of this run; remove it here (it is recreated at\n# its new home in Sprint 6 below).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\nReplace-InParagraph 7 \"Navigation and design\" \"Remaining Bug fixes + design (mobile first) + testing\"\n\n# Sprint 7: \"Navigation and design\" -> \"Navigation, Backbone History, design (mobile first) + testing\"\nReplace-InParagraph 6 \"Navigation and design\" \"Navigation, Backbone History, design (mobile first) + testing\"\n\n# Drop the \"_GoBack\" bookmark at the marker left in Sprint 6, then remove the marker itself.\n$markerRng = $d.Paragraphs(5).Range\n$markerFound = $markerRng.Find.Execute($marker)\nif (-not $markerFound) {\n    throw \"Could not find bookmark marker\"\n}\n$markerRng.Text = \"\"\n$d.Bookmarks.Add(\"_GoBack\", $markerRng)\n"}
